$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 1.30894078192739
$arr[0,1] = 0.142381227458003
$arr[0,2] = 0.4987193280759641
$arr[0,3] = 0.1610303018711168
$arr[0,4] = 0
$arr[0,5] = 0.8237303887048171
$arr[0,6] = 0.8931027468405262
$arr[0,7] = 0.8318112243451985
$arr[0,8] = 0.06234302277106618
$arr[0,9] = 0
$arr[0,10] = 0.4475299938882245
$arr[0,11] = 0
$arr[0,12] = 1.315702795095667
$arr[0,13] = 3.45520901001251
$arr[1,0] = 1.21536432012158
$arr[1,1] = 0.1310391777224424
$arr[1,2] = 0.4974918089285723
$arr[1,3] = 0.1618009621198855
$arr[1,4] = 0
$arr[1,5] = 0.8240526790278437
$arr[1,6] = 0.8975604967407804
$arr[1,7] = 0.8391401848351521
$arr[1,8] = 0.0625834658560569
$arr[1,9] = 0
$arr[1,10] = 0.4392220721310593
$arr[1,11] = 0
$arr[1,12] = 1.323399950736459
$arr[1,13] = 3.46491372699569
$arr[2,0] = 1.158110474741818
$arr[2,1] = 0.1240271722318056
$arr[2,2] = 0.4969543038967146
$arr[2,3] = 0.1623320022799071
$arr[2,4] = 0
$arr[2,5] = 0.8247716671149874
$arr[2,6] = 0.9006908013210904
$arr[2,7] = 0.8440747278029868
$arr[2,8] = 0.06274004911365427
$arr[2,9] = 0
$arr[2,10] = 0.4342931431209962
$arr[2,11] = 0
$arr[2,12] = 1.328600094259578
$arr[2,13] = 3.472789165687686
$arr[3,0] = 1.134831784498516
$arr[3,1] = 0.1211578432973681
$arr[3,2] = 0.4967897377116941
$arr[3,3] = 0.1625629751952449
$arr[3,4] = 0
$arr[3,5] = 0.8251956460268701
$arr[3,6] = 0.9020653813259827
$arr[3,7] = 0.8461948428458257
$arr[3,8] = 0.06280611444944029
$arr[3,9] = 0
$arr[3,10] = 0.4323280613176337
$arr[3,11] = 0
$arr[3,12] = 1.330838582953469
$arr[3,13] = 3.476480485872372
$arr[4,0] = 1.130969612732258
$arr[4,1] = 0.1206806810724288
$arr[4,2] = 0.4967657047159406
$arr[4,3] = 0.1626022087095471
$arr[4,4] = 0
$arr[4,5] = 0.8252739571566678
$arr[4,6] = 0.9022996081710062
$arr[4,7] = 0.8465534841728548
$arr[4,8] = 0.0628172209977933
$arr[4,9] = 0
$arr[4,10] = 0.4320043940764009
$arr[4,11] = 0
$arr[4,12] = 1.331217498469158
$arr[4,13] = 3.477122540645951
$arr[5,0] = 1.157796314348502
$arr[5,1] = 0.1239885233258207
$arr[5,2] = 0.4969518638039574
$arr[5,3] = 0.1623350582389982
$arr[5,4] = 0
$arr[5,5] = 0.8247768547543615
$arr[5,6] = 0.9007089386228415
$arr[5,7] = 0.844102878096237
$arr[5,8] = 0.06274093095036104
$arr[5,9] = 0
$arr[5,10] = 0.4342664649530832
$arr[5,11] = 0
$arr[5,12] = 1.328629799659794
$arr[5,13] = 3.472836996406897
$arr[6,0] = 1.276634642789759
$arr[6,1] = 0.1384805407912211
$arr[6,2] = 0.4982512705594075
$arr[6,3] = 0.1612840318030031
$arr[6,4] = 0
$arr[6,5] = 0.8237333278339918
$arr[6,6] = 0.8945582087124819
$arr[6,7] = 0.834248052127883
$arr[6,8] = 0.06242407340709732
$arr[6,9] = 0
$arr[6,10] = 0.4446298075267521
$arr[6,11] = 0
$arr[6,12] = 1.318258491028409
$arr[6,13] = 3.458157404352249
$arr[7,0] = 1.511215527179161
$arr[7,1] = 0.1665125085602881
$arr[7,2] = 0.5025107771809587
$arr[7,3] = 0.1596810544710934
$arr[7,4] = 0
$arr[7,5] = 0.8258250857491163
$arr[7,6] = 0.8856138468907773
$arr[7,7] = 0.8183715343836155
$arr[7,8] = 0.061873470918405
$arr[7,9] = 0
$arr[7,10] = 0.4663108349630818
$arr[7,11] = 0
$arr[7,12] = 1.301674233153889
$arr[7,13] = 3.444580489449208
$arr[8,0] = 1.684427446304142
$arr[8,1] = 0.1868649714924402
$arr[8,2] = 0.5066790581976903
$arr[8,3] = 0.1587813993752789
$arr[8,4] = 0
$arr[8,5] = 0.8298910260848089
$arr[8,6] = 0.8809394794111256
$arr[8,7] = 0.8088110560024759
$arr[8,8] = 0.06151171290063395
$arr[8,9] = 0
$arr[8,10] = 0.4830600405313987
$arr[8,11] = 0
$arr[8,12] = 1.291768595535061
$arr[8,13] = 3.443885034669051
$arr[9,0] = 1.763399525318789
$arr[9,1] = 0.1960698216262244
$arr[9,2] = 0.5088000012253815
$arr[9,3] = 0.1584322496639672
$arr[9,4] = 0
$arr[9,5] = 0.8322913886966461
$arr[9,6] = 0.8792242769614091
$arr[9,7] = 0.804918970918493
$arr[9,8] = 0.0613563496086007
$arr[9,9] = 0
$arr[9,10] = 0.4908561997362426
$arr[9,11] = 0
$arr[9,12] = 1.28775517128453
$arr[9,13] = 3.445585525547841
$arr[10,0] = 1.793328067524953
$arr[10,1] = 0.1995475986502697
$arr[10,2] = 0.5096353761234695
$arr[10,3] = 0.1583086581258826
$arr[10,4] = 0
$arr[10,5] = 0.8332796278190386
$arr[10,6] = 0.8786338436503058
$arr[10,7] = 0.8035108911582469
$arr[10,8] = 0.06129883496440591
$arr[10,9] = 0
$arr[10,10] = 0.4938336579255349
$arr[10,11] = 0
$arr[10,12] = 1.2863060912847
$arr[10,13] = 3.446519547921696
$arr[11,0] = 1.786881404796645
$arr[11,1] = 0.1987989515825745
$arr[11,2] = 0.509454031832675
$arr[11,3] = 0.1583348925365868
$arr[11,4] = 0
$arr[11,5] = 0.8330632658942108
$arr[11,6] = 0.8787583774058447
$arr[11,7] = 0.8038112204845689
$arr[11,8] = 0.06131116323313179
$arr[11,9] = 0
$arr[11,10] = 0.4931912906324101
$arr[11,11] = 0
$arr[11,12] = 1.286615033516171
$arr[11,13] = 3.446305486987228
$arr[12,0] = 1.765861304015516
$arr[12,1] = 0.1963560998738672
$arr[12,2] = 0.5088680828728656
$arr[12,3] = 0.1584219089901886
$arr[12,4] = 0
$arr[12,5] = 0.8323711025179819
$arr[12,6] = 0.8791745179716628
$arr[12,7] = 0.8048018089732878
$arr[12,8] = 0.06135159145766256
$arr[12,9] = 0
$arr[12,10] = 0.4911006527348576
$arr[12,11] = 0
$arr[12,12] = 1.287634538079161
$arr[12,13] = 3.445656554079278
$arr[13,0] = 1.752988892158953
$arr[13,1] = 0.1948587481027459
$arr[13,2] = 0.5085133651567588
$arr[13,3] = 0.1584763315847724
$arr[13,4] = 0
$arr[13,5] = 0.831957458542476
$arr[13,6] = 0.8794371080636267
$arr[13,7] = 0.80541713945928
$arr[13,8] = 0.06137652640951075
$arr[13,9] = 0
$arr[13,10] = 0.4898233550446633
$arr[13,11] = 0
$arr[13,12] = 1.288268219589241
$arr[13,13] = 3.445296843162538
$arr[14,0] = 1.679269833323076
$arr[14,1] = 0.1862623179670209
$arr[14,2] = 0.5065449638331359
$arr[14,3] = 0.1588054249810149
$arr[14,4] = 0
$arr[14,5] = 0.8297452457899936
$arr[14,6] = 0.8810598416745847
$arr[14,7] = 0.8090746135921698
$arr[14,8] = 0.06152205098355168
$arr[14,9] = 0
$arr[14,10] = 0.4825540840958951
$arr[14,11] = 0
$arr[14,12] = 1.292040785157404
$arr[14,13] = 3.443814493227507
$arr[15,0] = 1.634089559943106
$arr[15,1] = 0.180974823249187
$arr[15,2] = 0.5053949100824582
$arr[15,3] = 0.1590226942566915
$arr[15,4] = 0
$arr[15,5] = 0.8285292394883896
$arr[15,6] = 0.882160614847848
$arr[15,7] = 0.811435445199951
$arr[15,8] = 0.06161367879240398
$arr[15,9] = 0
$arr[15,10] = 0.4781397704367265
$arr[15,11] = 0
$arr[15,12] = 1.294481228131993
$arr[15,13] = 3.443421729462671
$arr[16,0] = 1.608119847544117
$arr[16,1] = 0.1779285638251622
$arr[16,2] = 0.5047545922065382
$arr[16,3] = 0.1591533206040605
$arr[16,4] = 0
$arr[16,5] = 0.827881659081072
$arr[16,6] = 0.8828324567765549
$arr[16,7] = 0.8128363474873943
$arr[16,8] = 0.06166724715259253
$arr[16,9] = 0
$arr[16,10] = 0.4756174336051515
$arr[16,11] = 0
$arr[16,12] = 1.29593128928947
$arr[16,13] = 3.443385654624052
$arr[17,0] = 1.599329891257867
$arr[17,1] = 0.1768962943455676
$arr[17,2] = 0.5045414293832096
$arr[17,3] = 0.1591985209431925
$arr[17,4] = 0
$arr[17,5] = 0.8276712996952682
$arr[17,6] = 0.8830665803499897
$arr[17,7] = 0.8133180552068531
$arr[17,8] = 0.06168553345491912
$arr[17,9] = 0
$arr[17,10] = 0.4747662814200311
$arr[17,11] = 0
$arr[17,12] = 1.296430225521057
$arr[17,13] = 3.443406042314024
$arr[18,0] = 1.638897353947414
$arr[18,1] = 0.1815382083222516
$arr[18,2] = 0.5055151457555667
$arr[18,3] = 0.1589989800481675
$arr[18,4] = 0
$arr[18,5] = 0.8286533202202264
$arr[18,6] = 0.8820394300837222
$arr[18,7] = 0.8111796784624445
$arr[18,8] = 0.06160383521763801
$arr[18,9] = 0
$arr[18,10] = 0.4786079589784151
$arr[18,11] = 0
$arr[18,12] = 1.294216639331331
$arr[18,13] = 3.443443892428434
$arr[19,0] = 1.772034793984574
$arr[19,1] = 0.197073840583613
$arr[19,2] = 0.509039316663376
$arr[19,3] = 0.1583961162543801
$arr[19,4] = 0
$arr[19,5] = 0.8325722555978103
$arr[19,6] = 0.8790506844909203
$arr[19,7] = 0.8045090637978234
$arr[19,8] = 0.06133968097796938
$arr[19,9] = 0
$arr[19,10] = 0.4917140409254728
$arr[19,11] = 0
$arr[19,12] = 1.287333166683268
$arr[19,13] = 3.445839288140405
$arr[20,0] = 1.859184060335963
$arr[20,1] = 0.2071811375457742
$arr[20,2] = 0.5115303203718469
$arr[20,3] = 0.1580523674529637
$arr[20,4] = 0
$arr[20,5] = 0.8355956105078235
$arr[20,6] = 0.8774416863115135
$arr[20,7] = 0.8005328026380951
$arr[20,8] = 0.061174721200981
$arr[20,9] = 0
$arr[20,10] = 0.5004265380034951
$arr[20,11] = 0
$arr[20,12] = 1.283246551724815
$arr[20,13] = 3.449095673225514
$arr[21,0] = 1.812659029017652
$arr[21,1] = 0.2017909685129098
$arr[21,2] = 0.5101836789979757
$arr[21,3] = 0.1582312405370914
$arr[21,4] = 0
$arr[21,5] = 0.8339396812577036
$arr[21,6] = 0.8782689507781072
$arr[21,7] = 0.8026199134497887
$arr[21,8] = 0.06126206230203479
$arr[21,9] = 0
$arr[21,10] = 0.4957631411651704
$arr[21,11] = 0
$arr[21,12] = 1.285389987660494
$arr[21,13] = 3.447202945817651
$arr[22,0] = 1.636723734463942
$arr[22,1] = 0.1812835218783562
$arr[22,2] = 0.5054607222146075
$arr[22,3] = 0.1590096834378425
$arr[22,4] = 0
$arr[22,5] = 0.8285970628523245
$arr[22,6] = 0.8820940962500998
$arr[22,7] = 0.8112951747065686
$arr[22,8] = 0.0616082827243325
$arr[22,9] = 0
$arr[22,10] = 0.478396242609108
$arr[22,11] = 0
$arr[22,12] = 1.294336113458783
$arr[22,13] = 3.443433281555258
$arr[23,0] = 1.447597957743824
$arr[23,1] = 0.1589712361452769
$arr[23,2] = 0.5011757163596116
$arr[23,3] = 0.1600657864848518
$arr[23,4] = 0
$arr[23,5] = 0.8248155299052655
$arr[23,6] = 0.8877001558636266
$arr[23,7] = 0.8222971802875101
$arr[23,8] = 0.06201488653722187
$arr[23,9] = 0
$arr[23,10] = 0.460300934017468
$arr[23,11] = 0
$arr[23,12] = 1.305759853083742
$arr[23,13] = 3.44662434687308

$ws.Range("B2:O25").Value = $arr